# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column (with
# constant value "stock") inserted right after the existing "total" column
# and before the existing "date" column. The columns that used to sit at
# H/I/J (date / legislator_name / legislator_id) therefore shift one column
# to the right, ending up at I/J/K.
#
# A few stock-name shared strings also had a stray embedded space removed
# as part of this same commit; fix those up while we're touching the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# ---- fix stray spaces inside a few company names (column B) ----
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $name = $cell.Value2
    if ($name -eq "★太平洋電線電纜股份有 限公司") {
        $cell.Value = "★太平洋電線電纜股份有限公司"
    } elseif ($name -eq "★春雨開發股份有限公司 (原正華）") {
        $cell.Value = "★春雨開發股份有限公司(原正華）"
    } elseif ($name -eq "★台灣上地開發股份有限 公司") {
        $cell.Value = "★台灣上地開發股份有限公司"
    }
}

# ---- shift H/I/J (date, legislator_name, legislator_id) one column right
#      to I/J/K, making room for the new H column ----
for ($r = 1; $r -le 13; $r++) {
    $dateCell = $ws.Cells.Item($r, 8)
    $legNameCell = $ws.Cells.Item($r, 9)
    $legIdCell = $ws.Cells.Item($r, 10)

    $dateVal = $dateCell.Value2
    $legNameVal = $legNameCell.Value2
    $legIdVal = $legIdCell.Value2

    $ws.Cells.Item($r, 11).Value = $legIdVal
    $ws.Cells.Item($r, 10).Value = $legNameVal
    $ws.Cells.Item($r, 9).Value = $dateVal
}

# ---- insert the new "property_category" column at H ----
$ws.Cells.Item(1, 8).Value = "property_category"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
